# Update cryptos list (GitHub Actions scheduled refresh) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.005.01"
$ws.Range("E2").Value = "  -1.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.174.34"
$ws.Range("E3").Value = "  -4.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "590.14"
$ws.Range("E5").Value = "  -2.63%  "

# Row 6 - Solana
$ws.Range("D6").Value = "134.33"
$ws.Range("E6").Value = "  -5.03%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.172.06"
$ws.Range("E8").Value = "  -4.01%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -1.31%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -6.17%  "

# Row 11 - Toncoin (unchanged)

# Row 12 - Cardano
$ws.Range("E12").Value = "  -3.60%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -4.73%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "34.24"
$ws.Range("E14").Value = "  -1.87%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.695.90"
$ws.Range("E15").Value = "  -4.06%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  -1.86%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.171.00"
$ws.Range("E17").Value = "  -4.11%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "62.985.38"
$ws.Range("E18").Value = "  -1.35%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -5.19%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "458.12"
$ws.Range("E20").Value = "  -4.61%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -0.45%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -6.20%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  -5.30%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "13.23"
$ws.Range("E24").Value = "  -5.28%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "82.26"
$ws.Range("E25").Value = "  -3.50%  "

# Row 26 - Dai
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.12%  "

# Row 27 - FirstDigitalUSD
$ws.Range("E27").Value = "  -0.11%  "

# Row 28 - PancakeSwap
$ws.Range("D28").Value = "2.65"
$ws.Range("E28").Value = "  -4.51%  "

# Row 29 - NEARProtocol
$ws.Range("E29").Value = "  -7.56%  "

# Row 30 - RenderToken
$ws.Range("D30").Value = "7.63"
$ws.Range("E30").Value = "  -6.80%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -6.04%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "27.14"
$ws.Range("E32").Value = "  -5.59%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -4.70%  "

# Row 34 - Stacks
$ws.Range("E34").Value = "  -6.81%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -6.36%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  -4.66%  "

# Row 37 - OKB
$ws.Range("D37").Value = "51.06"
$ws.Range("E37").Value = "  -2.68%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0704"
$ws.Range("E38").Value = "  -5.53%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -3.81%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "399.95"
$ws.Range("E40").Value = "  -7.18%  "

# Row 41 - Cosmos
$ws.Range("D41").Value = "8.08"
$ws.Range("E41").Value = "  -3.08%  "

# Row 42 - dogwifhat
$ws.Range("E42").Value = "  -4.00%  "

# Row 43 - was Maker, now Kaspa (rows 43/44 swapped order in the ranking)
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.111"
$ws.Range("E43").Value = "  -6.32%  "

# Row 44 - was Kaspa, now Maker
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.805.82"
$ws.Range("E44").Value = "  -9.87%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -5.52%  "

# Row 46 - USDe
$ws.Range("E46").Value = "  -0.02%  "

# Row 47 - Fetch.AI
$ws.Range("E47").Value = "  -5.74%  "

# Row 48 - Monero
$ws.Range("D48").Value = "124.41"
$ws.Range("E48").Value = "  -1.69%  "

# Row 49 - Arweave
$ws.Range("D49").Value = "34.54"
$ws.Range("E49").Value = "  -6.11%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "25.05"
$ws.Range("E50").Value = "  -5.02%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -2.49%  "
